$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.110.50"
$ws.Range("D2").Style = $ws.Range("B2").Style
$ws.Range("E2").Value = "  +5.36%  "
$ws.Range("D3").Value = "'1.880.98"
$ws.Range("D3").Style = $ws.Range("B3").Style
$ws.Range("E3").Value = "  +4.03%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'280.91"
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Value = "  +2.26%  "
$ws.Range("D6").Value = "'0.9995"
$ws.Range("D6").Style = $ws.Range("B6").Style
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "'0.5300"
$ws.Range("D7").Style = $ws.Range("B7").Style
$ws.Range("E7").Value = "  +5.09%  "
$ws.Range("D8").Value = "'0.3529"
$ws.Range("D8").Style = $ws.Range("B8").Style
$ws.Range("E8").Value = "  +0.57%  "
$ws.Range("D9").Value = "'45.36"
$ws.Range("D9").Style = $ws.Range("B9").Style
$ws.Range("E9").Value = "  +2.09%  "
$ws.Range("D10").Value = "'0.07067"
$ws.Range("D10").Style = $ws.Range("B10").Style
$ws.Range("E10").Value = "  +6.58%  "
$ws.Range("D11").Value = "'20.36"
$ws.Range("D11").Style = $ws.Range("B11").Style
$ws.Range("E11").Value = "  +2.00%  "
$ws.Range("D12").Value = "'0.8184"
$ws.Range("D12").Style = $ws.Range("B12").Style
$ws.Range("E12").Value = "  -1.61%  "
$ws.Range("D13").Value = "'0.07805"
$ws.Range("D13").Style = $ws.Range("B13").Style
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").Value = "'1.889.69"
$ws.Range("D14").Style = $ws.Range("B14").Style
$ws.Range("E14").Value = "  +4.56%  "
$ws.Range("D15").Value = "'5.207"
$ws.Range("D15").Style = $ws.Range("B15").Style
$ws.Range("E15").Value = "  +3.07%  "
$ws.Range("D16").Value = "'90.71"
$ws.Range("D16").Style = $ws.Range("B16").Style
$ws.Range("E16").Value = "  +3.92%  "
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "'14.59"
$ws.Range("D18").Style = $ws.Range("B18").Style
$ws.Range("E18").Value = "  +5.27%  "
$ws.Range("D19").Value = "'0.000008204"
$ws.Range("D19").Style = $ws.Range("B19").Style
$ws.Range("E19").Value = "  +2.82%  "
$ws.Range("D20").Value = "'0.9996"
$ws.Range("D20").Style = $ws.Range("B20").Style
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").Value = "'27.142.97"
$ws.Range("D21").Style = $ws.Range("B21").Style
$ws.Range("E21").Value = "  +5.19%  "
$ws.Range("D22").Value = "'2.124.54"
$ws.Range("D22").Style = $ws.Range("B22").Style
$ws.Range("E22").Value = "  +4.55%  "
$ws.Range("D23").Value = "'4.776"
$ws.Range("D23").Style = $ws.Range("B23").Style
$ws.Range("E23").Value = "  +1.25%  "
$ws.Range("E24").Value = "  +1.91%  "
$ws.Range("D25").Value = "'6.217"
$ws.Range("D25").Style = $ws.Range("B25").Style
$ws.Range("E25").Value = "  +3.01%  "
$ws.Range("D26").Value = "'2.388"
$ws.Range("D26").Style = $ws.Range("B26").Style
$ws.Range("E26").Value = "  +11.41%  "
$ws.Range("D27").Value = "'146.93"
$ws.Range("D27").Style = $ws.Range("B27").Style
$ws.Range("E27").Value = "  +3.89%  "
$ws.Range("D28").Value = "'17.60"
$ws.Range("D28").Style = $ws.Range("B28").Style
$ws.Range("E28").Value = "  +4.02%  "
$ws.Range("D29").Value = "'1.674"
$ws.Range("D29").Style = $ws.Range("B29").Style
$ws.Range("E29").Value = "  +1.21%  "
$ws.Range("D30").Value = "'112.74"
$ws.Range("D30").Style = $ws.Range("B30").Style
$ws.Range("E30").Value = "  +4.00%  "
$ws.Range("D31").Value = "'4.405"
$ws.Range("D31").Style = $ws.Range("B31").Style
$ws.Range("E31").Value = "  +1.85%  "
$ws.Range("D32").Value = "'4.392"
$ws.Range("D32").Style = $ws.Range("B32").Style
$ws.Range("E32").Value = "  +4.76%  "
$ws.Range("D33").Value = "'0.08937"
$ws.Range("D33").Style = $ws.Range("B33").Style
$ws.Range("E33").Value = "  +1.97%  "
$ws.Range("D34").Value = "'0.04910"
$ws.Range("D34").Style = $ws.Range("B34").Style
$ws.Range("E34").Value = "  +2.28%  "
$ws.Range("D35").Value = "'1.179"
$ws.Range("D35").Style = $ws.Range("B35").Style
$ws.Range("E35").Value = "  +4.15%  "
$ws.Range("D36").Value = "'0.7488"
$ws.Range("D36").Style = $ws.Range("B36").Style
$ws.Range("E36").Value = "  +3.18%  "
$ws.Range("D37").Value = "'2.906"
$ws.Range("D37").Style = $ws.Range("B37").Style
$ws.Range("E37").Value = "  +0.92%  "
$ws.Range("D38").Value = "'3.310"
$ws.Range("D38").Style = $ws.Range("B38").Style
$ws.Range("E38").Value = "  +9.07%  "
$ws.Range("D39").Value = "'2.420"
$ws.Range("D39").Style = $ws.Range("B39").Style
$ws.Range("E39").Value = "  +7.09%  "
$ws.Range("D40").Value = "'0.5339"
$ws.Range("D40").Style = $ws.Range("B40").Style
$ws.Range("E40").Value = "  +3.30%  "
$ws.Range("D41").Value = "'0.01884"
$ws.Range("D41").Style = $ws.Range("B41").Style
$ws.Range("E41").Value = "  +1.59%  "
$ws.Range("E42").Value = "  +3.08%  "
$ws.Range("D43").Value = "'117.14"
$ws.Range("D43").Style = $ws.Range("B43").Style
$ws.Range("E43").Value = "  +4.51%  "
$ws.Range("D44").Value = "'6.333"
$ws.Range("D44").Style = $ws.Range("B44").Style
$ws.Range("E44").Value = "  +2.89%  "
$ws.Range("D45").Value = "'8.238"
$ws.Range("D45").Style = $ws.Range("B45").Style
$ws.Range("E45").Value = "  +3.35%  "
$ws.Range("D46").Value = "'0.9987"
$ws.Range("D46").Style = $ws.Range("B46").Style
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").Value = "'0.4615"
$ws.Range("D47").Style = $ws.Range("B47").Style
$ws.Range("E47").Value = "  +1.55%  "
$ws.Range("D48").Value = "'0.1372"
$ws.Range("D48").Style = $ws.Range("B48").Style
$ws.Range("E48").Value = "  -0.40%  "
$ws.Range("D49").Value = "'9.486"
$ws.Range("D49").Style = $ws.Range("B49").Style
$ws.Range("E49").Value = "  +1.94%  "
$ws.Range("D50").Value = "'36.77"
$ws.Range("D50").Style = $ws.Range("B50").Style
$ws.Range("E50").Value = "  +2.25%  "
$ws.Range("D51").Value = "'1.533"
$ws.Range("D51").Style = $ws.Range("B51").Style
$ws.Range("E51").Value = "  +2.69%  "
